$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A5").Copy()
$ws.Range("A6").PasteSpecial(-4122)

$ws.Range("A6").Value = 44104
$ws.Range("B6").Value = "USD"
$ws.Range("C6").Value = "HKD"
$ws.Range("D6").Value = 7.75

$ws.Range("A7").Select()
